# Remove the <w:contextualSpacing w:val="0"/> element from every
# paragraph's pPr in the document (45 paragraphs in this file).
#
# The Word object model in this runtime does not expose
# ParagraphFormat.ContextualSpacing, so we can't toggle it through a
# dedicated property. Instead, for each paragraph we pull its current
# OOXML via Range.WordOpenXML, strip the <w:contextualSpacing/> element
# from that paragraph's fragment with a targeted string replace, and
# write the fragment back with Range.InsertXML so only that one element
# is removed - everything else in the paragraph (text, runs, other
# pPr children) is carried over unchanged.

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $full = $para.Range.WordOpenXML

    $bodyTag = "<w:body>"
    $bodyStart = $full.IndexOf($bodyTag)
    if ($bodyStart -lt 0) {
        continue
    }
    $searchFrom = $bodyStart + $bodyTag.Length

    $pStart = $full.IndexOf("<w:p", $searchFrom)
    if ($pStart -lt 0) {
        continue
    }
    $pEnd = $full.IndexOf("</w:p>", $pStart)
    if ($pEnd -lt 0) {
        continue
    }
    $pEnd = $pEnd + 6

    $fragment = $full.Substring($pStart, $pEnd - $pStart)

    if ($fragment -notmatch "<w:contextualSpacing\b[^/>]*/>") {
        continue
    }

    $newFragment = $fragment -replace "<w:contextualSpacing\b[^/>]*/>", ""

    [void]$para.Range.InsertXML($newFragment)
}

Write-Host "Removed contextualSpacing from $count paragraphs (where present)."
